# This script reproduces the "create core Dataset object" refactor for the
# merge_partial_expected_result.xlsx test fixture:
#   - the three "...link" lookup-key columns (A:C) are renamed to the
#     pandas merge-suffix style "..._x" (PIDN_x / DCDate_x / InstrID_x)
#   - the merge-diagnostic columns (K:M) are renamed from the old
#     "_merge" / "_diff_days" / "_abs_diff_days" names to the new
#     "_mp_merge" / "_mp_diff_days" / "_mp_abs_diff_days" names
#   - the trailing "_duplicates" column (N), which was always FALSE, is
#     removed entirely
#   - column widths for the renamed K:M columns are widened to fit the
#     new, longer header text

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the merge-diagnostic header columns (values/position unchanged) ---
$ws.Range("K1").Value = "_mp_merge"
$ws.Range("L1").Value = "_mp_diff_days"
$ws.Range("M1").Value = "_mp_abs_diff_days"

# --- Rename the lookup-key header columns (values/position unchanged) ---
$ws.Range("A1").Value = "PIDN_x"
$ws.Range("B1").Value = "DCDate_x"
$ws.Range("C1").Value = "InstrID_x"

# --- Drop the old "_duplicates" column (N) entirely, shifting nothing else ---
$ws.Columns.Item(14).Delete()

# --- Widen K:L:M to bestFit the longer renamed headers ---
# (Excel's ColumnWidth is in character units and gets rounded to whole
# pixels on save, so we target the character width whose rounded pixel
# width matches the saved file as closely as possible.)
$ws.Columns.Item(11).ColumnWidth = 9.498697916666666
$ws.Columns.Item(12).ColumnWidth = 11.166666666666666
$ws.Columns.Item(13).ColumnWidth = 14.498697916666666
